$wb = $excel.ActiveWorkbook

# The "4000" series sheet is renamed to "CMOS" and becomes the active
# (selected) sheet/tab, with cell C18 selected on it.
$ws = $wb.Worksheets.Item("4000")
$ws.Name = "CMOS"
$ws.Activate()
$ws.Range("C18").Select()
